# The log analyzer was re-run after 3 more log lines were appended to the source log:
# 3 additional ServiceB ERROR entries ('Job X failed to start') at 09:12:01, 09:13:01 and
# 09:44:01. Refresh the derived sheets accordingly and add a new summary sheet that
# counts how often each distinct error message occurs.

$wb = $excel.ActiveWorkbook

# ---- 1. Successful Entries: append the 3 newly parsed log rows ----
$wsSuccess = $wb.Worksheets.Item("Successful Entries")
$wsSuccess.Range("A13").Value = '2023-03-01 09:12:01'
$wsSuccess.Range("B13").Value = 'ServiceB'
$wsSuccess.Range("C13").Value = 'ERROR'
$wsSuccess.Range("D13").Value = 'Job X failed to start'
$wsSuccess.Range("A14").Value = '2023-03-01 09:13:01'
$wsSuccess.Range("B14").Value = 'ServiceB'
$wsSuccess.Range("C14").Value = 'ERROR'
$wsSuccess.Range("D14").Value = 'Job X failed to start'
$wsSuccess.Range("A15").Value = '2023-03-01 09:44:01'
$wsSuccess.Range("B15").Value = 'ServiceB'
$wsSuccess.Range("C15").Value = 'ERROR'
$wsSuccess.Range("D15").Value = 'Job X failed to start'

# ---- 2. Service Entries: refresh ServiceB's count + entries list ----
$wsService = $wb.Worksheets.Item("Service Entries")
$wsService.Range("B3").Value = 7
$wsService.Range("C3").Value = '[{''date_time'': ''2023-03-01 08:15:28'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Null pointer exception''}, {''date_time'': ''2023-03-01 08:35:10'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Null pointer exception''}, {''date_time'': ''2023-03-01 09:00:00'', ''service_name'': ''ServiceB'', ''log_level'': ''INFO'', ''message'': ''Started job X''}, {''date_time'': ''2023-03-01 09:00:01'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Job X failed to start''}, {''date_time'': ''2023-03-01 09:12:01'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Job X failed to start''}, {''date_time'': ''2023-03-01 09:13:01'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Job X failed to start''}, {''date_time'': ''2023-03-01 09:44:01'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Job X failed to start''}]'

# ---- 3. Log Level Entries: refresh ERROR's count + entries list ----
$wsLogLevel = $wb.Worksheets.Item("Log Level Entries")
$wsLogLevel.Range("B3").Value = 6
$wsLogLevel.Range("C3").Value = '[{''date_time'': ''2023-03-01 08:15:28'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Null pointer exception''}, {''date_time'': ''2023-03-01 08:35:10'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Null pointer exception''}, {''date_time'': ''2023-03-01 09:00:01'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Job X failed to start''}, {''date_time'': ''2023-03-01 09:12:01'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Job X failed to start''}, {''date_time'': ''2023-03-01 09:13:01'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Job X failed to start''}, {''date_time'': ''2023-03-01 09:44:01'', ''service_name'': ''ServiceB'', ''log_level'': ''ERROR'', ''message'': ''Job X failed to start''}]'

# ---- 4. New sheet: Error Messages count (grouped count of ERROR messages) ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsErrCount = $wb.Worksheets.Add($null, $lastSheet)
$wsErrCount.Name = " Error Messages count "

# Reuse the existing bold/bordered header formatting from the other summary sheets
# instead of rebuilding it cell-by-cell.
$headerSrc = $wsLogLevel.Range("A1:B1")
$headerSrc.Copy()
$wsErrCount.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsErrCount.Range("A1").Value = "Error Message"
$wsErrCount.Range("B1").Value = "Count"
$wsErrCount.Range("A2").Value = "Null pointer exception"
$wsErrCount.Range("B2").Value = 2
$wsErrCount.Range("A3").Value = "Job X failed to start"
$wsErrCount.Range("B3").Value = 4

Write-Output "edit complete"